$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 96.666664
$ws.Range("I2").Value = 96.666664
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 96.666664
$ws.Range("L2").Value = 0
$ws.Range("M2").ClearContents()
$ws.Range("N2").Value = 16.333336

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 333
$ws.Range("I18").Value = 299.5
$ws.Range("J18").Value = 400
$ws.Range("K18").Value = 299.5
$ws.Range("L18").Value = 400
$ws.Range("M18").Value = -15.5
$ws.Range("N18").Value = -968

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 1171.6666
$ws.Range("I58").Value = 257.5
$ws.Range("J58").Value = 3000
$ws.Range("K58").Value = 772.5
$ws.Range("L58").Value = 9000
$ws.Range("M58").Value = -622.5
$ws.Range("N58").Value = -9300

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 25001178
$ws.Range("I137").Value = 37037890
$ws.Range("J137").Value = 1858.5385
$ws.Range("K137").Value = 111113670
$ws.Range("L137").Value = 5575.6155
$ws.Range("M137").Value = -111111120
$ws.Range("N137").Value = -10675.6155

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 11473.546
$ws.Range("I86").Value = 3740.4
$ws.Range("J86").Value = 17917.834
$ws.Range("K86").Value = 3740.4
$ws.Range("L86").Value = 17917.834
$ws.Range("M86").Value = -2617.4
$ws.Range("N86").Value = -20163.834

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 11473.546
$ws.Range("I89").Value = 3740.4
$ws.Range("J89").Value = 17917.834
$ws.Range("K89").Value = 18702
$ws.Range("L89").Value = 89589.17
$ws.Range("M89").Value = -13086
$ws.Range("N89").Value = -100821.17

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1548.36
$ws.Range("I99").Value = 1456.3182
$ws.Range("K99").Value = 1456.3182
$ws.Range("M99").Value = 41.68180000000007

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 866.85
$ws.Range("I107").Value = 825.64703
$ws.Range("K107").Value = 825.64703
$ws.Range("M107").Value = 1094.35297

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1991.3948
$ws.Range("I31").Value = 1128.2413
$ws.Range("K31").Value = 1128.2413
$ws.Range("M31").Value = -833.2412999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 1991.3948
$ws.Range("I34").Value = 1128.2413
$ws.Range("K34").Value = 1128.2413
$ws.Range("M34").Value = -926.2412999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 12300
$ws.Range("I50").Value = 5000
$ws.Range("J50").Value = 14733.333
$ws.Range("K50").Value = 5000
$ws.Range("L50").Value = 14733.333
$ws.Range("M50").Value = -4375
$ws.Range("N50").Value = -15983.333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H51").Value = 19000
$ws.Range("I51").Value = 18000
$ws.Range("J51").Value = 21000
$ws.Range("K51").Value = 18000
$ws.Range("L51").Value = 21000
$ws.Range("M51").Value = -17264
$ws.Range("N51").Value = -22472

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").ClearContents()
$ws.Range("N52").Value = 0

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H59").Value = 20000
$ws.Range("J59").Value = 20000
$ws.Range("L59").Value = 20000
$ws.Range("N59").Value = -22290

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H61").Value = 19000
$ws.Range("I61").Value = 18000
$ws.Range("J61").Value = 21000
$ws.Range("K61").Value = 18000
$ws.Range("L61").Value = 21000
$ws.Range("M61").Value = -17652
$ws.Range("N61").Value = -21696

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 24849.5
$ws.Range("I62").Value = 26499.445
$ws.Range("J62").Value = 10000
$ws.Range("K62").Value = 26499.445
$ws.Range("L62").Value = 10000
$ws.Range("M62").Value = -25875.445
$ws.Range("N62").Value = -11248

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 24849.5
$ws.Range("I65").Value = 26499.445
$ws.Range("J65").Value = 10000
$ws.Range("K65").Value = 132497.225
$ws.Range("L65").Value = 50000
$ws.Range("M65").Value = -129377.225
$ws.Range("N65").Value = -56240

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H68").Value = 30000
$ws.Range("J68").Value = 30000
$ws.Range("L68").Value = 30000
$ws.Range("N68").Value = -31498

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H71").Value = 30000
$ws.Range("J71").Value = 30000
$ws.Range("L71").Value = 90000
$ws.Range("N71").Value = -97488

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H74").Value = 32966.668
$ws.Range("J74").Value = 44950
$ws.Range("L74").Value = 44950
$ws.Range("N74").Value = -46698

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H77").Value = 32966.668
$ws.Range("J77").Value = 44950
$ws.Range("L77").Value = 134850
$ws.Range("N77").Value = -143586

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H133").Value = 22374
$ws.Range("I133").Value = 20296
$ws.Range("J133").Value = 22670.857
$ws.Range("K133").Value = 20296
$ws.Range("L133").Value = 22670.857
$ws.Range("M133").Value = -17766
$ws.Range("N133").Value = -27730.857

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H135").Value = 40814.4
$ws.Range("J135").Value = 40814.4
$ws.Range("L135").Value = 40814.4
$ws.Range("N135").Value = -50954.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 2214.7368
$ws.Range("I113").Value = 1930.375
$ws.Range("K113").Value = 1930.375
$ws.Range("M113").Value = 239.625

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3062.3635
$ws.Range("I132").Value = 2727.4866
$ws.Range("J132").Value = 4832.4287
$ws.Range("K132").Value = 8182.459800000001
$ws.Range("L132").Value = 14497.2861
$ws.Range("M132").Value = -5652.459800000001
$ws.Range("N132").Value = -19557.2861

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2295.5557
$ws.Range("I46").Value = 1475
$ws.Range("J46").Value = 2952
$ws.Range("K46").Value = 1475
$ws.Range("L46").Value = 2952
$ws.Range("M46").Value = -1287
$ws.Range("N46").Value = -3328

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 17857.715
$ws.Range("I61").Value = 17834
$ws.Range("J61").Value = 18000
$ws.Range("K61").Value = 17834
$ws.Range("L61").Value = 18000
$ws.Range("M61").Value = -17632
$ws.Range("N61").Value = -18404

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2520
$ws.Range("J68").Value = 3033.3333
$ws.Range("L68").Value = 3033.3333
$ws.Range("N68").Value = -4531.3333

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 2520
$ws.Range("J71").Value = 3033.3333
$ws.Range("L71").Value = 15166.6665
$ws.Range("N71").Value = -22654.6665

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 17857.715
$ws.Range("I113").Value = 17834
$ws.Range("J113").Value = 18000
$ws.Range("K113").Value = 17834
$ws.Range("L113").Value = 18000
$ws.Range("M113").Value = -15664
$ws.Range("N113").Value = -22340

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 3304.6545
$ws.Range("I136").Value = 2191.2083
$ws.Range("J136").Value = 10939.714
$ws.Range("K136").Value = 6573.624899999999
$ws.Range("L136").Value = 32819.142
$ws.Range("M136").Value = -4023.624899999999
$ws.Range("N136").Value = -37919.142

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").ClearContents()
$ws.Range("N110").Value = 0
